# "Continued work on part 2" - mark Part 2's Finished date in the
# timetable table, and move the sheet's active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Part 2 is the table row where Part = 1 (B4), Deadline = 19.06.2022 (C4).
# Its "Finished" column cell (D4) had no value yet - fill in the date the
# part was actually finished.
$tbl = $ws.ListObjects.Item("Table1")
$finishedCol = $tbl.ListColumns.Item("Finished")
$finishedCol.DataBodyRange.Cells.Item(2, 1).Value = "21.06.2022"

# Move the selection up one row, to F7.
$ws.Range("F7").Select()
